$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell C1 - reuse the existing header style (same as A1/B1)
$ws.Range("C1").Value = "Complemento"
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats

# New (empty) cell C2 alongside the existing data row 2
$ws.Range("C2").Font.Bold = $false

# New data row 3
$ws.Range("A3").Value = "Intimissimi"
$ws.Range("B3").Value = "Av"
$ws.Range("C3").Font.Bold = $false
